# Updates cryptos list values (price / volume change %) and restores the
# original Cosmos/Toncoin row order swap, per the upstream data refresh.
# Cell styling is preserved: NumberFormat is temporarily forced to Text
# while assigning the value (so numeric-looking strings like "419.73" or
# "0.663" are not coerced into floating point numbers), then the original
# Style is restored so no new style entries end up applied to the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '66.203.23'
$cell.Style = $origStyle

$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.48%  '
$cell.Style = $origStyle

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.550.96'
$cell.Style = $origStyle

$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.77%  '
$cell.Style = $origStyle

$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell.Style = $origStyle

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '419.73'
$cell.Style = $origStyle

$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.33%  '
$cell.Style = $origStyle

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '132.17'
$cell.Style = $origStyle

$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.59%  '
$cell.Style = $origStyle

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.663'
$cell.Style = $origStyle

$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.71%  '
$cell.Style = $origStyle

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.542.48'
$cell.Style = $origStyle

$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.74%  '
$cell.Style = $origStyle

$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell.Style = $origStyle

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.784'
$cell.Style = $origStyle

$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +8.40%  '
$cell.Style = $origStyle

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.168'
$cell.Style = $origStyle

$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +20.85%  '
$cell.Style = $origStyle

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000288'
$cell.Style = $origStyle

$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +34.14%  '
$cell.Style = $origStyle

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '43.43'
$cell.Style = $origStyle

$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.87%  '
$cell.Style = $origStyle

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.10'
$cell.Style = $origStyle

$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +8.76%  '
$cell.Style = $origStyle

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.122.38'
$cell.Style = $origStyle

$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.77%  '
$cell.Style = $origStyle

$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.24%  '
$cell.Style = $origStyle

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '20.55'
$cell.Style = $origStyle

$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.50%  '
$cell.Style = $origStyle

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.560.41'
$cell.Style = $origStyle

$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.86%  '
$cell.Style = $origStyle

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.88'
$cell.Style = $origStyle

$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.95%  '
$cell.Style = $origStyle

$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.57%  '
$cell.Style = $origStyle

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '66.150.50'
$cell.Style = $origStyle

$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.33%  '
$cell.Style = $origStyle

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '449.66'
$cell.Style = $origStyle

$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.57%  '
$cell.Style = $origStyle

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '90.48'
$cell.Style = $origStyle

$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.27%  '
$cell.Style = $origStyle

$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.88%  '
$cell.Style = $origStyle

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.24'
$cell.Style = $origStyle

$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.58%  '
$cell.Style = $origStyle

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.39'
$cell.Style = $origStyle

$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.38%  '
$cell.Style = $origStyle

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.04'
$cell.Style = $origStyle

$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.83%  '
$cell.Style = $origStyle

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '34.27'
$cell.Style = $origStyle

$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.79%  '
$cell.Style = $origStyle

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.83'
$cell.Style = $origStyle

$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.05%  '
$cell.Style = $origStyle

$cell = $ws.Range("B30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Toncoin'
$cell.Style = $origStyle

$cell = $ws.Range("C30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell.Style = $origStyle

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.80'
$cell.Style = $origStyle

$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.89%  '
$cell.Style = $origStyle

$cell = $ws.Range("B31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Cosmos'
$cell.Style = $origStyle

$cell = $ws.Range("C31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell.Style = $origStyle

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.52'
$cell.Style = $origStyle

$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.91%  '
$cell.Style = $origStyle

$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.74%  '
$cell.Style = $origStyle

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.31'
$cell.Style = $origStyle

$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.29%  '
$cell.Style = $origStyle

$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.71%  '
$cell.Style = $origStyle

$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.06%  '
$cell.Style = $origStyle

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '39.29'
$cell.Style = $origStyle

$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.19%  '
$cell.Style = $origStyle

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '57.67'
$cell.Style = $origStyle

$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.40%  '
$cell.Style = $origStyle

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0508'
$cell.Style = $origStyle

$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.44%  '
$cell.Style = $origStyle

$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +42.38%  '
$cell.Style = $origStyle

$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +11.42%  '
$cell.Style = $origStyle

$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.15%  '
$cell.Style = $origStyle

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.04'
$cell.Style = $origStyle

$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.15%  '
$cell.Style = $origStyle

$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.64%  '
$cell.Style = $origStyle

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.46'
$cell.Style = $origStyle

$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.91%  '
$cell.Style = $origStyle

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '147.09'
$cell.Style = $origStyle

$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.18%  '
$cell.Style = $origStyle

$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.54%  '
$cell.Style = $origStyle

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.311'
$cell.Style = $origStyle

$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.46%  '
$cell.Style = $origStyle

$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.46%  '
$cell.Style = $origStyle

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell.Style = $origStyle

$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.01%  '
$cell.Style = $origStyle

$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.41%  '
$cell.Style = $origStyle

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.83'
$cell.Style = $origStyle

$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.45%  '
$cell.Style = $origStyle
